# Emission plots and postprocess:
# The default emission name "CO2-equivalent" is replaced with a user
# specified emission name ("CO2") plus a separate descriptive name
# ("CO2 emissions"), adding a new "emission_name" column to the
# Emissions sheet (which pushes "emission_unit" / "ton" one column right).

$wb = $excel.ActiveWorkbook
$ws4 = $wb.Worksheets.Item("Emissions")

# --- Restructure the Emissions sheet from 2 columns (A,B) to 3 (A,B,C) ---

# Give the new header cells (A1, C1) the same header style (bold/border/
# center) already used by B1, by copying formats only.
$ws4.Range("B1").Copy()
$ws4.Range("A1").PasteSpecial(-4122)
$ws4.Range("C1").PasteSpecial(-4122)

# Move the old "emission_unit" header from B1 -> C1, and put the new
# "emission_name" header in B1.
$ws4.Range("C1").Value = $ws4.Range("B1").Value()
$ws4.Range("B1").Value = "emission_name"

# Move the old unit value ("ton") from B2 -> C2, replace the old
# "CO2-equivalent" single value (A2) with the short code "CO2", and add
# the longer descriptive name in the new B2 cell.
$ws4.Range("C2").Value = $ws4.Range("B2").Value()
$ws4.Range("A2").Value = "CO2"
$ws4.Range("B2").Value = "CO2 emissions"

# Resize the columns to fit the new content.
$ws4.Columns("A:B").ColumnWidth = 12.3
$ws4.Columns("C:C").ColumnWidth = 11.15

# --- Misc selection updates on the other (non-active) sheets ---

$ws1 = $wb.Worksheets.Item("Techs")
$ws1.Range("A1:B2").Select()

$ws2 = $wb.Worksheets.Item("Fuels")
$ws2.Range("E12").Select()

# Finally, make Emissions the active (tabSelected) sheet and set its
# selection / active cell. This must be last so Emissions ends up as the
# workbook's active tab.
$ws4.Activate()
$ws4.Range("E9").Select()
